$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.971.72"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.640.72"
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "'216.13"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "1.870.35"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "1.644.02"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "'62.88"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "25.934.98"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'193.04"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +7.39%  "
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "'144.54"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'3.28"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "1.134.45"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "'99.45"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "1.779.74"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("D46").Value = "'56.65"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'7.76"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  +0.73%  "
